$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "2025-08-03 06:42:05"
$ws.Range("B11").Value = "delete-team"
$ws.Range("C11").Value = "new-organization97"
$ws.Range("D11").Value = "firstteam"
$ws.Range("I11").Value = "'False"
$ws.Range("I11").Style = "Normal"
